# Portal_Verification.xlsx - the "Verifications" sheet is a flat, single
# column (A) list of question/label strings, one per row, used to drive an
# automated test. This adds two new rows for the new "Bypass URL checks"
# test case: a heading row and the appeal-reference row it verifies.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A171").Value = "Exemption validation"
$ws.Range("A172").Value = "Appeal a planning decision / 848-HAS"
